# Add a new "2022-Q3" quarterly sheet to the workbook.
#
# 1. Insert a new worksheet named "2022-Q3" right before the current
#    "2022-Q2" sheet (so tab order becomes 总计, 2022-Q3, 2022-Q2, ...).
# 2. Populate it with the fund holdings table for that quarter.
# 3. Insert a new row into "总计" (summary) sheet for 2022-Q3 totals,
#    pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing "store as text" semantics
# (mirrors the workbook's existing convention of keeping numeric-looking
# figures such as "97.93" as inline strings instead of numbers).
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Step 1: create + position the new "2022-Q3" sheet
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# The other quarterly sheets share one layout: a bold/bordered header row
# (B1:H1) and a bold/bordered running-index column (A2:A{n}). Borrow that
# formatting from "2022-Q2" (same column layout) so the new sheet matches.
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2:A10").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: header row
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $newSheet.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# Step 3: fund holding rows (code, name, scale, stockPosition, positionPct, marketValue, rank)
# ---------------------------------------------------------------------------
$funds = @(
    @("516950", "银华中证基建ETF", "11.07", "97.93", "3.74", "0.4140", 8),
    @("159635", "华夏中证基建ETF", "3.40", "99.03", "3.73", "0.1268", 8),
    @("159619", "国泰中证基建ETF", "3.30", "98.76", "3.69", "0.1218", 8),
    @("000646", "华润元大量化优选混合A", "1.47", "73.62", "4.59", "0.0675", 7),
    @("004260", "德邦稳盈增长灵活配置混合", "1.37", "88.81", "4.87", "0.0667", 5),
    @("001412", "德邦鑫星价值灵活配置混合A", "0.97", "87.27", "4.83", "0.0469", 6),
    @("007827", "华润元大量化优选混合C", "0.19", "73.62", "4.59", "0.0087", 7),
    @("002112", "德邦鑫星价值灵活配置混合C", "0.08", "87.27", "4.83", "0.0039", 6),
    @("519165", "新华鑫利灵活配置混合", "0.05", "25.03", "3.38", "0.0017", 3)
)

$r = 2
$idx = 0
foreach ($fund in $funds) {
    $newSheet.Cells.Item($r, 1).Value = $idx

    Set-TextCell $newSheet.Cells.Item($r, 2) $fund[0]
    Set-TextCell $newSheet.Cells.Item($r, 3) $fund[1]
    Set-TextCell $newSheet.Cells.Item($r, 4) $fund[2]
    Set-TextCell $newSheet.Cells.Item($r, 5) $fund[3]
    Set-TextCell $newSheet.Cells.Item($r, 6) $fund[4]
    Set-TextCell $newSheet.Cells.Item($r, 7) $fund[5]

    $newSheet.Cells.Item($r, 8).Value = $fund[6]

    $r = $r + 1
    $idx = $idx + 1
}

# ---------------------------------------------------------------------------
# Step 4: update the "总计" (summary) sheet with a new 2022-Q3 row
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows("2:2").Insert()

# The Insert() shifted the old row-2..row-8 data down to row-3..row-9, but
# left the old row-2 (now empty) without formatting, and kept the shifted
# rows' column-A running index untouched (still 0..6 instead of 1..7).
# Re-apply the data-row formatting to the new row 2 by copying it from row 3
# (which still carries the original per-column styling) ...
$summary.Range("B3:D3").Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# ... then fill in the new 2022-Q3 figures ...
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 0.86

# ... and renumber the running index (column A) of every shifted row so it
# keeps counting 0,1,2,3,... down the table.
for ($row = 3; $row -le 9; $row++) {
    $summary.Cells.Item($row, 1).Value = $row - 2
}
